$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the old hyperlinks (rows 2-15) before we touch row layout ---
$ws.Range("F2:F15").Hyperlinks.Delete()

# --- Drop the now-stale rows 8-15 (the refreshed scrape only keeps 6 listings) ---
$ws.Rows.Item(8).Resize(8).EntireRow.Delete()

# --- Widen column H (13 -> 17 chars). COM ColumnWidth has a fixed +0.8333 padding
#     offset versus the stored OOXML <col width>, so back it out to land on 17. ---
$ws.Columns.Item(8).ColumnWidth = 16.166666666666666

# --- Refreshed listing data (2026-01-28 06:33:29 JST scrape) ---
$timestamp = "2026-01-28 06:33:29"

$ws.Cells.Item(2,1).Value = $timestamp
$ws.Cells.Item(2,2).Value = "Difyと連携したAIチャットアプリ(Bubble)の実証実験用プロダクト(MVP)の開発"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5480502"
$ws.Cells.Item(2,7).Value = 378
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◆開発 ◇アプリ"

$ws.Cells.Item(3,1).Value = $timestamp
$ws.Cells.Item(3,2).Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5455098"
$ws.Cells.Item(3,7).Value = 375
$ws.Cells.Item(3,8).Value = "🔥AI,Ai ◆開発"

$ws.Cells.Item(4,1).Value = $timestamp
$ws.Cells.Item(4,2).Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5445159"
$ws.Cells.Item(4,7).Value = 368
$ws.Cells.Item(4,8).Value = "🔥AI,Ai ◆開発"

$ws.Cells.Item(5,1).Value = $timestamp
$ws.Cells.Item(5,2).Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5445154"
$ws.Cells.Item(5,7).Value = 368
$ws.Cells.Item(5,8).Value = "🔥AI,Ai ◆開発"

$ws.Cells.Item(6,1).Value = $timestamp
$ws.Cells.Item(6,2).Value = "Instagram投稿をWordPressへ自動連携したいです(公式API/将来拡張あり)"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5480464"
$ws.Cells.Item(6,7).Value = 208
$ws.Cells.Item(6,8).Value = "🔥API ○WordPress"

$ws.Cells.Item(7,1).Value = $timestamp
$ws.Cells.Item(7,2).Value = "急募 【1.6万円/即決】超シンプル仕様の3分タイマーアプリ3本(iOS/Android計6ビルド)の開発"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "10,000 円 ~ 20,000 円 / 募集期間 3 日、取引期間 0 日"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5480298"
$ws.Cells.Item(7,7).Value = 85
$ws.Cells.Item(7,8).Value = "◆開発 ◇アプリ"

# --- Re-create hyperlinks for the URL column (F2:F7) ---
$urls = @{
    2 = "https://www.lancers.jp/work/detail/5480502"
    3 = "https://www.lancers.jp/work/detail/5455098"
    4 = "https://www.lancers.jp/work/detail/5445159"
    5 = "https://www.lancers.jp/work/detail/5445154"
    6 = "https://www.lancers.jp/work/detail/5480464"
    7 = "https://www.lancers.jp/work/detail/5480298"
}
foreach ($r in 2..7) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$r])
}
